# Minor changes to documentation:
# The "Algorithms and Complexity" skill entry (row 9, column C) is corrected
# to "Algorithms and Complexity Analysis".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Algorithms and Complexity Analysis"

# Leave the active selection on the edited cell, matching the saved state.
$ws.Range("C9").Select()
